$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (28/07/2022): Invalid and Absent marked
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4 (01/08/2022): Total Attendance Count and Real marked
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5 (04/08/2022): Total Attendance Count and Real marked
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Rows 6-18: Absent marked for each remaining date
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
